{"js": "// Office.js (Word JavaScript API) edit script.\n// The document consists of a single lead paragraph (a date string) followed\n// by a 20x5 table of arithmetic equations (100 cells, one paragraph each).\n// Each paragraph in document order holds exactly one run of text that needs\n// to be replaced with a new value (old -> new pairs taken from the diff, in\n// document order).\n\nconst REPLACEMENTS = [\n  [\"2026-02-19 Thursday\", \"2026-02-20 Friday\"],\n  [\"71-18=\", \"71+1=\"],\n  [\"3+83=\", \"85-13=\"],\n  [\"71-44=\", \"80-58=\"],\n  [\"28-17=\", \"88-28=\"],\n  [\"95-10=\", \"31-17=\"],\n  [\"7+76=\", \"80+4=\"],\n  [\"99-47=\", \"23+68=\"],\n  [\"91-2=\", \"49+8=\"],\n  [\"1+19=\", \"94-5=\"],\n  [\"14+16=\", \"77-58=\"],\n  [\"97-11=\", \"52+30=\"],\n  [\"95-23=\", \"96-44=\"],\n  [\"91+0=\", \"70-47=\"],\n  [\"17+77=\", \"78-18=\"],\n  [\"21-9=\", \"42-16=\"],\n  [\"51-34=\", \"13+13=\"],\n  [\"96-77=\", \"54-7=\"],\n  [\"79+11=\", \"70-38=\"],\n  [\"94-59=\", \"77-69=\"],\n  [\"94-4=\", \"38+44=\"],\n  [\"62-43=\", \"73-40=\"],\n  [\"65-64=\", \"37+53=\"],\n  [\"3+31=\", \"6+88=\"],\n  [\"87-10=\", \"46+51=\"],\n  [\"75-31=\", \"84+7=\"],\n  [\"38+49=\", \"20-17=\"],\n  [\"94-8=\", \"86-55=\"],\n  [\"85-24=\", \"67-41=\"],\n  [\"94-89=\", \"17+71=\"],\n  [\"22+56=\", \"12+61=\"],\n  [\"68-50=\", \"4+23=\"],\n  [\"74-71=\", \"0+84=\"],\n  [\"55-10=\", \"57+13=\"],\n  [\"18+10=\", \"75-2=\"],\n  [\"69-23=\", \"70-41=\"],\n  [\"75-27=\", \"13-5=\"],\n  [\"26+19=\", \"66-52=\"],\n  [\"31-8=\", \"27-22=\"],\n  [\"73-68=\", \"19+0=\"],\n  [\"83-23=\", \"69+8=\"],\n  [\"34+48=\", \"89-53=\"],\n  [\"92-69=\", \"43+56=\"],\n  [\"66+9=\", \"59-13=\"],\n  [\"82-78=\", \"59-7=\"],\n  [\"66-20=\", \"17+59=\"],\n  [\"21+54=\", \"65-55=\"],\n  [\"73-67=\", \"32+6=\"],\n  [\"69-68=\", \"14+1=\"],\n  [\"75-47=\", \"26+25=\"],\n  [\"98-57=\", \"27+1=\"],\n  [\"83-65=\", \"30+64=\"],\n  [\"1+81=\", \"24+22=\"],\n  [\"71-47=\", \"26+69=\"],\n  [\"85-84=\", \"9+41=\"],\n  [\"27-26=\", \"29+56=\"],\n  [\"52-35=\", \"6+38=\"],\n  [\"75-53=\", \"97-4=\"],\n  [\"91-50=\", \"25-7=\"],\n  [\"0+79=\", \"13+80=\"],\n  [\"50+2=\", \"53+19=\"],\n  [\"33+60=\", \"83-66=\"],\n  [\"93-36=\", \"25+16=\"],\n  [\"42+39=\", \"26+66=\"],\n  [\"88+3=\", \"51+23=\"],\n  [\"11+2=\", \"52+32=\"],\n  [\"34-29=\", \"95-41=\"],\n  [\"33+6=\", \"74-27=\"],\n  [\"44+42=\", \"98-87=\"],\n  [\"82-32=\", \"98-30=\"],\n  [\"9+31=\", \"15+40=\"],\n  [\"3+56=\", \"4+48=\"],\n  [\"62-48=\", \"56+8=\"],\n  [\"72-37=\", \"11+69=\"],\n  [\"28+39=\", \"79+7=\"],\n  [\"80-17=\", \"16+77=\"],\n  [\"17+5=\", \"40-7=\"],\n  [\"21+74=\", \"5+1=\"],\n  [\"97-76=\", \"89-32=\"],\n  [\"4+31=\", \"26+65=\"],\n  [\"14+5=\", \"87-54=\"],\n  [\"84-28=\", \"50+49=\"],\n  [\"45+18=\", \"53+3=\"],\n  [\"24+29=\", \"29+28=\"],\n  [\"64-40=\", \"11+79=\"],\n  [\"52-25=\", \"25+73=\"],\n  [\"41-40=\", \"31+59=\"],\n  [\"0+81=\", \"57-24=\"],\n  [\"85-70=\", \"46+49=\"],\n  [\"56+14=\", \"11+22=\"],\n  [\"42+11=\", \"90-87=\"],\n  [\"14+46=\", \"59-24=\"],\n  [\"19+46=\", \"23+62=\"],\n  [\"45+14=\", \"92-16=\"],\n  [\"89-1=\", \"31-28=\"],\n  [\"33+51=\", \"2+25=\"],\n  [\"1+17=\", \"4+64=\"],\n  [\"43-9=\", \"49-37=\"],\n  [\"29-24=\", \"5+45=\"],\n  [\"13+32=\", \"21+25=\"],\n  [\"23-0=\", \"83-27=\"],\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    \"Unexpected paragraph count: found \" + items.length +\n    \", expected \" + REPLACEMENTS.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  // Defensive check: confirm we are editing the paragraph the diff expects.\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: found \" + JSON.stringify(para.text) +\n      \", expected \" + JSON.stringify(oldText)\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# PowerShell / Word COM-interop edit script.\n# Document = 1 lead paragraph (date) + a 20x5 table of arithmetic equations.\n# old/new text pairs below come from the diff, in document order.\n\n$DateOld = '2026-02-19 Thursday'\n$DateNew = '2026-02-20 Friday'\n\n# Row-major (row 1..20, each with 5 columns) old/new equation text.\n$CellOld = @(\n    '71-18=', '3+83=', '71-44=', '28-17=', '95-10=',\n    '7+76=', '99-47=', '91-2=', '1+19=', '14+16=',\n    '97-11=', '95-23=', '91+0=', '17+77=', '21-9=',\n    '51-34=', '96-77=', '79+11=', '94-59=', '94-4=',\n    '62-43=', '65-64=', '3+31=', '87-10=', '75-31=',\n    '38+49=', '94-8=', '85-24=', '94-89=', '22+56=',\n    '68-50=', '74-71=', '55-10=', '18+10=', '69-23=',\n    '75-27=', '26+19=', '31-8=', '73-68=', '83-23=',\n    '34+48=', '92-69=', '66+9=', '82-78=', '66-20=',\n    '21+54=', '73-67=', '69-68=', '75-47=', '98-57=',\n    '83-65=', '1+81=', '71-47=', '85-84=', '27-26=',\n    '52-35=', '75-53=', '91-50=', '0+79=', '50+2=',\n    '33+60=', '93-36=', '42+39=', '88+3=', '11+2=',\n    '34-29=', '33+6=', '44+42=', '82-32=', '9+31=',\n    '3+56=', '62-48=', '72-37=', '28+39=', '80-17=',\n    '17+5=', '21+74=', '97-76=', '4+31=', '14+5=',\n    '84-28=', '45+18=', '24+29=', '64-40=', '52-25=',\n    '41-40=', '0+81=', '85-70=', '56+14=', '42+11=',\n    '14+46=', '19+46=', '45+14=', '89-1=', '33+51=',\n    '1+17=', '43-9=', '29-24=', '13+32=', '23-0='\n)\n\n$CellNew = @(\n    '71+1=', '85-13=', '80-58=', '88-28=', '31-17=',\n    '80+4=', '23+68=', '49+8=', '94-5=', '77-58=',\n    '52+30=', '96-44=', '70-47=', '78-18=', '42-16=',\n    '13+13=', '54-7=', '70-38=', '77-69=', '38+44=',\n    '73-40=', '37+53=', '6+88=', '46+51=', '84+7=',\n    '20-17=', '86-55=', '67-41=', '17+71=', '12+61=',\n    '4+23=', '0+84=', '57+13=', '75-2=', '70-41=',\n    '13-5=', '66-52=', '27-22=', '19+0=', '69+8=',\n    '89-53=', '43+56=', '59-13=', '59-7=', '17+59=',\n    '65-55=', '32+6=', '14+1=', '26+25=', '27+1=',\n    '30+64=', '24+22=', '26+69=', '9+41=', '29+56=',\n    '6+38=', '97-4=', '25-7=', '13+80=', '53+19=',\n    '83-66=', '25+16=', '26+66=', '51+23=', '52+32=',\n    '95-41=', '74-27=', '98-87=', '98-30=', '15+40=',\n    '4+48=', '56+8=', '11+69=', '79+7=', '16+77=',\n    '40-7=', '5+1=', '89-32=', '26+65=', '87-54=',\n    '50+49=', '53+3=', '29+28=', '11+79=', '25+73=',\n    '31+59=', '57-24=', '46+49=', '11+22=', '90-87=',\n    '59-24=', '23+62=', '92-16=', '31-28=', '2+25=',\n    '4+64=', '49-37=', '5+45=', '21+25=', '83-27='\n)\n\n$d = $word.ActiveDocument\n\n# Word's Range.Text includes trailing paragraph (\\r) / cell (\\a) marks;\n# strip those before comparing against plain expected strings.\nfunction Strip-Marks($text) {\n    return $text.TrimEnd([char]13, [char]7)\n}\n\n# --- Update the lead date paragraph ---\n$dateParagraph = $d.Paragraphs.Item(1)\n$currentDateText = Strip-Marks $dateParagraph.Range.Text\nif ($currentDateText -ne $DateOld) {\n    throw \"Lead paragraph text mismatch: found '$currentDateText', expected '$DateOld'\"\n}\n$dateParagraph.Range.Text = $DateNew\n\n# --- Update the 20x5 table of equations ---\n$table = $d.Tables.Item(1)\nif ($table.Rows.Count -ne 20 -or $table.Columns.Count -ne 5) {\n    throw \"Unexpected table size: $($table.Rows.Count) rows x $($table.Columns.Count) cols\"\n}\n\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $table.Cell($r, $c)\n        $expectedOld = $CellOld[($r - 1) * 5 + ($c - 1)]\n        $expectedNew = $CellNew[($r - 1) * 5 + ($c - 1)]\n        $currentText = Strip-Marks $cell.Range.Text\n        if ($currentText -ne $expectedOld) {\n            throw \"Cell ($r,$c) text mismatch: found '$currentText', expected '$expectedOld'\"\n        }\n        if ($expectedOld -ne $expectedNew) {\n            $cell.Range.Text = $expectedNew\n        }\n    }\n}\n"}
